$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 900.4167
$ws.Range("I98").Value = 422.77777
$ws.Range("J98").Value = 2333.3333
$ws.Range("K98").Value = 422.77777
$ws.Range("L98").Value = 2333.3333
$ws.Range("M98").Value = 1075.22223
$ws.Range("N98").Value = -5329.3333
$ws.Range("H122").Value = 900.4167
$ws.Range("I122").Value = 422.77777
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 1268.33331
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = 1181.66669
$ws.Range("N122").Value = -11899.9999
$ws.Range("H137").Value = 3624.0881
$ws.Range("I137").Value = 2597.16
$ws.Range("J137").Value = 6476.6665
$ws.Range("K137").Value = 7791.48
$ws.Range("L137").Value = 19429.9995
$ws.Range("M137").Value = -5241.48
$ws.Range("N137").Value = -24529.9995
$ws.Range("H138").Value = 2138.1267
$ws.Range("I138").Value = 2227
$ws.Range("J138").Value = 2105.6538
$ws.Range("K138").Value = 6681
$ws.Range("L138").Value = 6316.9614
$ws.Range("M138").Value = -1541
$ws.Range("N138").Value = -16596.9614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 430446.6
$ws.Range("I32").Value = 492517.25
$ws.Range("K32").Value = 492517.25
$ws.Range("M32").Value = -492230.25
$ws.Range("H61").Value = 4001.4375
$ws.Range("I61").Value = 4003.8333
$ws.Range("K61").Value = 4003.8333
$ws.Range("M61").Value = -3791.8333
$ws.Range("H74").Value = 2684.7
$ws.Range("I74").Value = 2369.4
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 2369.4
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -1495.4
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 2684.7
$ws.Range("I77").Value = 2369.4
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 11847
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -7479
$ws.Range("N77").Value = -23736
$ws.Range("H132").Value = 2698.1836
$ws.Range("I132").Value = 2416.1316
$ws.Range("J132").Value = 3672.5454
$ws.Range("K132").Value = 7248.3948
$ws.Range("L132").Value = 11017.6362
$ws.Range("M132").Value = -4718.3948
$ws.Range("N132").Value = -16077.6362
$ws.Range("H136").Value = 4001.4375
$ws.Range("I136").Value = 4003.8333
$ws.Range("K136").Value = 12011.4999
$ws.Range("M136").Value = -9461.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3114.4546
$ws.Range("I134").Value = 2962.111
$ws.Range("K134").Value = 8886.332999999999
$ws.Range("M134").Value = -6351.332999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6758.45
$ws.Range("I31").Value = 1473.1111
$ws.Range("J31").Value = 17735.691
$ws.Range("K31").Value = 1473.1111
$ws.Range("L31").Value = 17735.691
$ws.Range("M31").Value = -1178.1111
$ws.Range("N31").Value = -18325.691
$ws.Range("H34").Value = 6758.45
$ws.Range("I34").Value = 1473.1111
$ws.Range("J34").Value = 17735.691
$ws.Range("K34").Value = 1473.1111
$ws.Range("L34").Value = 17735.691
$ws.Range("M34").Value = -1271.1111
$ws.Range("N34").Value = -18139.691
$ws.Range("H58").Value = 1056.909
$ws.Range("I58").Value = 736
$ws.Range("J58").Value = 1550.6154
$ws.Range("K58").Value = 736
$ws.Range("L58").Value = 1550.6154
$ws.Range("M58").Value = -533
$ws.Range("N58").Value = -1956.6154
$ws.Range("H109").Value = 37166.668
$ws.Range("J109").Value = 37166.668
$ws.Range("L109").Value = 37166.668
$ws.Range("N109").Value = -39246.668
$ws.Range("H132").Value = 11113282
$ws.Range("I132").Value = 1450.6666
$ws.Range("J132").Value = 27781028
$ws.Range("K132").Value = 4351.9998
$ws.Range("L132").Value = 83343084
$ws.Range("M132").Value = -1821.9998
$ws.Range("N132").Value = -83348144
$ws.Range("H134").Value = 7601
$ws.Range("I134").Value = 6000
$ws.Range("J134").Value = 8668.333000000001
$ws.Range("K134").Value = 18000
$ws.Range("L134").Value = 26004.999
$ws.Range("M134").Value = -15465
$ws.Range("N134").Value = -31074.999
$ws.Range("H136").Value = 1056.909
$ws.Range("I136").Value = 736
$ws.Range("J136").Value = 1550.6154
$ws.Range("K136").Value = 2208
$ws.Range("L136").Value = 4651.8462
$ws.Range("M136").Value = 342
$ws.Range("N136").Value = -9751.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 2509.6775
$ws.Range("I64").Value = 650
$ws.Range("J64").Value = 2637.9312
$ws.Range("K64").Value = 1950
$ws.Range("L64").Value = 7913.7936
$ws.Range("M64").Value = -1680
$ws.Range("N64").Value = -8453.793600000001
$ws.Range("H67").Value = 2509.6775
$ws.Range("I67").Value = 650
$ws.Range("J67").Value = 2637.9312
$ws.Range("K67").Value = 1950
$ws.Range("L67").Value = 7913.7936
$ws.Range("M67").Value = -1014
$ws.Range("N67").Value = -9785.793600000001
$ws.Range("H74").Value = 2500
$ws.Range("J74").Value = 4000
$ws.Range("L74").Value = 12000
$ws.Range("N74").Value = -14122
$ws.Range("H77").Value = 2500
$ws.Range("J77").Value = 4000
$ws.Range("L77").Value = 36000
$ws.Range("N77").Value = -46608
$ws.Range("H81").Value = 7003
$ws.Range("J81").Value = 8999.666999999999
$ws.Range("L81").Value = 26999.001
$ws.Range("N81").Value = -29245.001
$ws.Range("H84").Value = 7003
$ws.Range("J84").Value = 8999.666999999999
$ws.Range("L84").Value = 80997.003
$ws.Range("N84").Value = -92229.003
$ws.Range("H110").Value = 9997.036
$ws.Range("I110").Value = 5506.75
$ws.Range("J110").Value = 10745.417
$ws.Range("K110").Value = 16520.25
$ws.Range("L110").Value = 32236.251
$ws.Range("M110").Value = -12430.25
$ws.Range("N110").Value = -40416.251
$ws.Range("H139").Value = 3446
$ws.Range("I139").Value = 3066.125
$ws.Range("J139").Value = 3749.9
$ws.Range("K139").Value = 9198.375
$ws.Range("L139").Value = 11249.7
$ws.Range("M139").Value = -4058.375
$ws.Range("N139").Value = -21529.7
$ws.Range("H141").Value = 4792
$ws.Range("I141").Value = 1345.3846
$ws.Range("J141").Value = 7150.2104
$ws.Range("K141").Value = 4036.1538
$ws.Range("L141").Value = 21450.6312
$ws.Range("M141").Value = 1143.8462
$ws.Range("N141").Value = -31810.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 786567.0600000001
$ws.Range("I80").Value = 2254501.2
$ws.Range("J80").Value = 52600
$ws.Range("K80").Value = 2254501.2
$ws.Range("L80").Value = 52600
$ws.Range("M80").Value = -2253503.2
$ws.Range("N80").Value = -54596
$ws.Range("H83").Value = 786567.0600000001
$ws.Range("I83").Value = 2254501.2
$ws.Range("J83").Value = 52600
$ws.Range("K83").Value = 11272506
$ws.Range("L83").Value = 263000
$ws.Range("M83").Value = -11267514
$ws.Range("N83").Value = -272984
$ws.Range("H132").Value = 2729
$ws.Range("I132").Value = 2240.9375
$ws.Range("J132").Value = 3844.5715
$ws.Range("K132").Value = 6722.8125
$ws.Range("L132").Value = 11533.7145
$ws.Range("M132").Value = -4192.8125
$ws.Range("N132").Value = -16593.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4812.3
$ws.Range("I132").Value = 4021.5
$ws.Range("J132").Value = 5998.5
$ws.Range("K132").Value = 12064.5
$ws.Range("L132").Value = 17995.5
$ws.Range("M132").Value = -9534.5
$ws.Range("N132").Value = -23055.5
$ws.Range("H134").Value = 39115.668
$ws.Range("J134").Value = 39115.668
$ws.Range("L134").Value = 39115.668
$ws.Range("N134").Value = -49255.668
$ws.Range("H136").Value = 7247993.5
$ws.Range("I136").Value = 1179.875
$ws.Range("K136").Value = 3539.625
$ws.Range("M136").Value = -989.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 4016000
$ws.Range("I109").Value = 8000000
$ws.Range("J109").Value = 32000
$ws.Range("K109").Value = 8000000
$ws.Range("L109").Value = 32000
$ws.Range("M109").Value = -7998613
$ws.Range("N109").Value = -34774
$ws.Range("H132").Value = 4169058
$ws.Range("I132").Value = 2949.5334
$ws.Range("J132").Value = 6668723
$ws.Range("K132").Value = 8848.600199999999
$ws.Range("L132").Value = 20006169
$ws.Range("M132").Value = -6318.600199999999
$ws.Range("N132").Value = -20011229
$ws.Range("H136").Value = 3330.6316
$ws.Range("I136").Value = 2892.625
$ws.Range("J136").Value = 5666.6665
$ws.Range("K136").Value = 8677.875
$ws.Range("L136").Value = 16999.9995
$ws.Range("M136").Value = -6127.875
$ws.Range("N136").Value = -22099.9995
